$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage (avoid Excel auto-converting numeric-looking strings) for the price/volume columns
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "26.983.73"
$ws.Range("E2").Value = "  -0.11%  "

$ws.Range("D3").Value = "1.844.00"
$ws.Range("E3").Value = "  -0.05%  "

$ws.Range("E4").Value = "  +0.58%  "

$ws.Range("D5").Value = "1.014"
$ws.Range("E5").Value = "  +0.54%  "

$ws.Range("D6").Value = "309.10"
$ws.Range("E6").Value = "  -0.40%  "

$ws.Range("D7").Value = "0.4755"
$ws.Range("E7").Value = "  +1.93%  "

$ws.Range("D8").Value = "0.3675"
$ws.Range("E8").Value = "  +1.30%  "

$ws.Range("D9").Value = "0.07217"
$ws.Range("E9").Value = "  +0.93%  "

$ws.Range("D10").Value = "0.9299"
$ws.Range("E10").Value = "  +1.45%  "

$ws.Range("D11").Value = "19.81"
$ws.Range("E11").Value = "  +1.02%  "

$ws.Range("D12").Value = "0.07758"
$ws.Range("E12").Value = "  +0.71%  "

$ws.Range("D13").Value = "1.842.55"
$ws.Range("E13").Value = "  -1.19%  "

$ws.Range("D14").Value = "5.386"
$ws.Range("E14").Value = "  +1.91%  "

$ws.Range("D15").Value = "6.465"
$ws.Range("E15").Value = "  +0.78%  "

$ws.Range("D16").Value = "88.77"
$ws.Range("E16").Value = "  +0.38%  "

$ws.Range("E17").Value = "  +0.59%  "

$ws.Range("D18").Value = "0.000008657"
$ws.Range("E18").Value = "  +0.81%  "

$ws.Range("E19").Value = "  +0.56%  "

$ws.Range("D20").Value = "27.025.12"
$ws.Range("E20").Value = "  -0.04%  "

$ws.Range("D21").Value = "14.52"
$ws.Range("E21").Value = "  +0.86%  "

$ws.Range("D22").Value = "5.054"
$ws.Range("E22").Value = "  +0.49%  "

$ws.Range("D23").Value = "10.63"
$ws.Range("E23").Value = "  -0.08%  "

$ws.Range("E24").Value = "  -0.01%  "

$ws.Range("D25").Value = "152.74"
$ws.Range("E25").Value = "  -0.02%  "

$ws.Range("D26").Value = "18.23"
$ws.Range("E26").Value = "  -0.67%  "

$ws.Range("D27").Value = "1.988"
$ws.Range("E27").Value = "  -3.73%  "

$ws.Range("D28").Value = "114.48"
$ws.Range("E28").Value = "  +0.22%  "

$ws.Range("D29").Value = "4.949"
$ws.Range("E29").Value = "  +0.61%  "

$ws.Range("D30").Value = "0.08862"
$ws.Range("E30").Value = "  +0.01%  "

$ws.Range("D31").Value = "3.307"
$ws.Range("E31").Value = "  +3.77%  "

$ws.Range("E32").Value = "  +0.40%  "

$ws.Range("E33").Value = "  +0.72%  "

$ws.Range("E34").Value = "  -1.45%  "

$ws.Range("D35").Value = "2.669"
$ws.Range("E35").Value = "  -6.92%  "

$ws.Range("D36").Value = "1.111"
$ws.Range("E36").Value = "  +2.56%  "

$ws.Range("D37").Value = "0.01966"
$ws.Range("E37").Value = "  +1.28%  "

$ws.Range("D38").Value = "0.05261"

$ws.Range("D39").Value = "2.974"
$ws.Range("E39").Value = "  -0.24%  "

$ws.Range("D40").Value = "0.5237"
$ws.Range("E40").Value = "  +1.24%  "

$ws.Range("D41").Value = "7.026"
$ws.Range("E41").Value = "  +1.79%  "

$ws.Range("D42").Value = "0.1511"
$ws.Range("E42").Value = "  +0.02%  "

$ws.Range("D43").Value = "8.270"

$ws.Range("D44").Value = "10.51"
$ws.Range("E44").Value = "  -0.10%  "

$ws.Range("D45").Value = "0.4729"
$ws.Range("E45").Value = "  +0.59%  "

$ws.Range("E46").Value = "  +0.53%  "

$ws.Range("D47").Value = "101.68"
$ws.Range("E47").Value = "  +0.90%  "

$ws.Range("D48").Value = "1.607"
$ws.Range("E48").Value = "  +0.06%  "

$ws.Range("D49").Value = "65.55"
$ws.Range("E49").Value = "  +0.98%  "

$ws.Range("D50").Value = "0.06061"
$ws.Range("E50").Value = "  +0.29%  "

$ws.Range("D51").Value = "0.8907"

# Restore default (General) formatting/style so cells match original unstyled cells
$ws.Range("D2:E51").ClearFormats()
